# ---------------------------------------------------------------------------
# Finish analyzing KI67 manual-vs-Aperio: build a "Percentage" summary row on
# Sheet2 (header row of case numbers + a row of percentages copied from
# Sheet1), and a sorted two-block case/score table on Sheet3. Sheet3 becomes
# the active sheet (as last worked on).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Source data lifted from Sheet1 row 1 (case ids) and row 18 (percentages),
# column by column B..N.
$cols  = @("B","C","D","E","F","G","H","I","J","K","L","M","N")
$cases = @(46, 74, 141, 152, 188, 226, 300, 316, 319, 326, 689, 622, 668)
$pcts  = @(47.890625, 15.949428640894725, 85.714285714285708, 29.745042492917843, `
           7.6180482686253939, 25.354609929078016, 46.268656716417908, `
           60.975609756097562, 11.760966306420851, 15.051020408163266, `
           72.941176470588232, 100, 22.388059701492537)

# =============================== Sheet2 ====================================
# Row 1: case-id header, same formatting (bold / bold+yellow) as Sheet1 row 1.
# Row 2: "Percentage" label + the raw percentage values (plain formatting,
# pasted as values only - matches Sheet1 row 18's numbers without its
# formulas/format).

$ws1.Range("A1:N1").Copy()
$ws2.Range("A1:N1").PasteSpecial(-4122)

$ws2.Range("A1").Value = "Case"
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws2.Range($cols[$i] + "1").Value = $cases[$i]
}

$ws2.Range("A2").Value = "Percentage"
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws2.Range($cols[$i] + "2").Value = $pcts[$i]
}

# =============================== Sheet3 ====================================
# Columns A:B - every case/percentage pair (same order as Sheet1 columns
# B..N), formatting per-cell copied from the matching Sheet1 row-1 cell.
# Columns D:E - the subset of cases that used Sheet1's plain-bold style
# (as opposed to bold+yellow), with their case id under header "case" and
# percentage under header "score".

$ws3.Range("A1").Value = "Case"
$ws1.Range("A1").Copy()
$ws3.Range("A1").PasteSpecial(-4122)
$ws3.Range("B1").Value = "Percentage"

$ws3.Range("D1").Value = "case"
$ws1.Range("A1").Copy()
$ws3.Range("D1").PasteSpecial(-4122)
$ws3.Range("E1").Value = "score"

for ($i = 0; $i -lt $cols.Length; $i++) {
    $row = $i + 2
    $ws1.Range($cols[$i] + "1").Copy()
    $ws3.Range("A" + $row).PasteSpecial(-4122)
    $ws3.Range("A" + $row).Value = $cases[$i]
    $ws3.Range("B" + $row).Value = $pcts[$i]
}

$subsetCols  = @("C","F","J","M","N")
$subsetCase  = @(74, 188, 319, 622, 668)
$subsetPct   = @(15.949428640894725, 7.6180482686253939, 11.760966306420851, 100, 22.388059701492537)

for ($i = 0; $i -lt $subsetCols.Length; $i++) {
    $row = $i + 2
    $ws1.Range($subsetCols[$i] + "1").Copy()
    $ws3.Range("D" + $row).PasteSpecial(-4122)
    $ws3.Range("D" + $row).Value = $subsetCase[$i]
    $ws3.Range("E" + $row).Value = $subsetPct[$i]
}

# =============================== Views/selection ============================
# Sheet1: active cell moves to the percentage row (row 18), selected as a
# whole row like a user clicking the row header.
$ws1.Rows.Item(18).Select()

# Sheet2: the pasted block A1:N2 is selected.
$ws2.Activate()
$ws2.Range("A1:N2").Select()

# Sheet3: becomes the active sheet, with the "case/score" block selected.
$ws3.Activate()
$ws3.Range("D1:E6").Select()
